$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header text in B1 (was "الرقم", now "رقم الهوية")
$ws.Range("B1").Value = "رقم الهوية"

# Move the active selection to B2 (as in the saved file)
$ws.Range("B2").Select()
